# The original "Sheet" had stray/invalid data (including a bogus row "0"
# that Excel cannot address/clear through the normal object model) left
# over from a prior run. The fix recreates the worksheet from scratch -
# i.e. the ".xlsx sheet is automatically (re)created if it does not
# exist / is not in a good state" - and then writes the single valid
# data row back onto it.

$wb = $excel.ActiveWorkbook

# Remember the name of the sheet we need to replace.
$sheetName = $wb.ActiveSheet.Name

# Add a brand-new, empty worksheet to stand in for it.
$freshSheet = $wb.Worksheets.Add()
$stagingName = "NewSheetStaging"
$freshSheet.Name = $stagingName

# Remove the old sheet (re-resolved by name, since index-bound
# references shift when a new sheet is inserted).
$wb.Worksheets.Item($sheetName).Delete() | Out-Null

# Rename the fresh sheet back to the original name so the workbook's
# sheet list/order/name is unchanged.
$ws = $wb.Worksheets.Item($stagingName)
$ws.Name = $sheetName

# Populate the (only) data row with the new subscriber entry. The id
# columns are digit strings (not numbers), so force text formatting
# before assigning them - otherwise Excel would auto-convert the
# numeric-looking text into a number.
$ws.Cells.Item(1, 1).Value = "ve"
$ws.Cells.Item(1, 2).Value = "dds"
$ws.Cells.Item(1, 3).NumberFormat = "@"
$ws.Cells.Item(1, 3).Value = "12345678"
$ws.Cells.Item(1, 4).NumberFormat = "@"
$ws.Cells.Item(1, 4).Value = "12345678"
$ws.Cells.Item(1, 5).Value = "GN"
$ws.Cells.Item(1, 6).Value = "النحل"
